$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 was corrupted: the reservation id ("asdf") had been typed into the
# guest_name-shaped slot as text in column A, and the actual guest name was
# never entered (column B empty). Fix the data:
#   A4 -> numeric reservation id (matches the pattern of rows 2-3: 12345, 12346)
#   B4 -> the guest's name
$ws.Range("A4").Value = 12348
$ws.Range("B4").Value = "Jakiś Pan"

# A4 already carried the row's normal cell formatting; toggling WrapText off
# (its already-default state) forces Excel to materialize that formatting as
# an explicit style so it survives cleanly now that the cell's content type
# changed from text to a number.
$ws.Range("A4").WrapText = $false

# Copy that same formatting onto the newly-populated B4 so both edited cells
# in the row share one consistent style, matching the rest of the table.
$ws.Range("A4").Copy() | Out-Null
$ws.Range("B4").PasteSpecial(-4122) | Out-Null
